$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.114.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.06%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.626.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.73%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.77%  '

# Row 10
$ws.Range("E10").Value = '  -0.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.332'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.78%  '

# Row 12
$ws.Range("E12").Value = '  +0.53%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.095.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.89%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.091.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.56'
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.641.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.31%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.09%  '

# Row 22
$ws.Range("E22").Value = '  +0.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.411'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.74%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.756.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.66%  '

# Row 27
$ws.Range("E27").Value = '  +0.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.39%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0780'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.98%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.74%  '

# Row 32
$ws.Range("E32").Value = '  +0.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.88'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.42%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.93%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.17%  '

# Row 38
$ws.Range("E38").Value = '  -7.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.834'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.68%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.91%  '

# Row 42
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0970'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.12%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.594'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.13%  '

# Row 45
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.75%  '

# Row 46
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '267.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.95%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.81%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0528'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.54%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.023.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0227'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.12%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.64%  '
